$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: MagneticSensor - fill in the "claimed" info
$ws.Range("F44").Value = "Magnetic"
$ws.Range("B44").Value = "HiTechnicMagneticSensor"
$ws.Range("D44").Value = "Lawrie"
$ws.Range("E44").Value = "N"
$ws.Range("G44").Value = "SampleProvider"

# Row 61: RCXLightSensor - fill in the "claimed" info
$ws.Range("F61").Value = "Light"
$ws.Range("B61").Value = "RCXLightSensor"
$ws.Range("D61").Value = "Lawrie"
$ws.Range("E61").Value = "N"
$ws.Range("G61").Value = "SampleProvider"

# Rows 62-68: mark "Fits in framework" column C
$ws.Range("C62").Value = "N"
$ws.Range("C63").Value = "N"
$ws.Range("C64").Value = "N"
$ws.Range("C65").Value = "N"
$ws.Range("C66").Value = "N"
$ws.Range("C67").Value = "?"
$ws.Range("C68").Value = "N"

# Match the author's final selection position
$ws.Range("C68").Select()
